$p = $ppt.ActivePresentation

# Slide 13: "Multivariate Analysis – Part 2" -> "Bivariate Analysis – Part 3"
$s13 = $p.Slides.Item(13)
$s13.Shapes.Item(1).TextFrame.TextRange.Text = "Bivariate Analysis – Part 3"

# Slide 14: "Bivariate Analysis – Part 3" -> "Bivariate Analysis – Part 4"
$s14 = $p.Slides.Item(14)
$s14.Shapes.Item(1).TextFrame.TextRange.Text = "Bivariate Analysis – Part 4"

# Slide 15: "Recommendation" -> "Recommendations"
$s15 = $p.Slides.Item(15)
$s15.Shapes.Item(1).TextFrame.TextRange.Text = "Recommendations"
